$d = $word.ActiveDocument

# --- Paragraph 1: update the placeholder id text, drop the trailing space run ---
$p1 = $d.Paragraphs(1)
$full = $p1.Range
$idRange = $d.Range($full.Start, $full.End - 2)
$idRange.Text = "**ID__AFFARS_SUBPART_5322_3__ID**"

$p1 = $d.Paragraphs(1)
$full = $p1.Range
$spaceRange = $d.Range($full.End - 2, $full.End - 1)
$spaceRange.Text = ""

# --- Paragraph 1: paragraph formatting (indent + border) ---
$p1 = $d.Paragraphs(1)
$pf = $p1.Range.ParagraphFormat
$pf.LeftIndent = 11.25
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5
